# Apply cryptos list update (coinranking data refresh) to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D (Price) and E (Volume(1h)) columns are stored as text in the source data;
# force Text number format before assigning so values such as "307.98" or
# "  +2.48%  " are not auto-coerced into numeric cells by COM.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "44.173.10"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.48%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.425.18"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.13%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.98"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.62%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.84"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +3.98%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.52%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.501"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.49%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.27"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +3.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0799"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.82%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "18.90"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +3.14%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.08%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.86%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.804.55"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.441.75"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.04%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.836"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +3.49%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "44.129.57"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.25"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.54%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.40"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.76%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0905"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.58"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.34%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +5.38%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "240.44"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.26%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.85%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.26"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.75%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.62%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.57"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +4.77%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.80"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +4.82%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +11.25%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.63"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +6.99%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.18%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0762"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.75%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +3.11%  "
$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "130.17"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +24.46%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.47"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +4.39%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +3.87%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.07%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.91%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.30"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -4.54%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.36%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.950.54"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.37%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +2.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.87"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +4.73%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.43"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +3.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.65"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +9.50%  "
$ws.Range("B49").Value = "MultiversX"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "53.43"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.53%  "
$ws.Range("B50").Value = "BitcoinSV"
$ws.Range("C50").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "73.78"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.57%  "
$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.15"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.78%  "
